$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 22: "SD" label + SD value (Base Gas Fees standard deviation) ---
# Set the "SD" shared string first so it is allocated before "SD EUR"
# (keeps shared-string table ordering identical to the authored workbook).
$ws.Range("A22").Value = "SD"
$ws.Range("B22").Value = 47409562748.343384

# --- New column K: "SD EUR" header + per-row formulas ---
# Copy J3's format (bold header style) onto K3 before setting its text.
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("K3").Value = "SD EUR"

# Match number formatting used by the rest of the numeric columns (two
# decimal places, same as F:J).
$ws.Range("K4:K7").NumberFormat = "0.00"

$ws.Range("K4").Formula = "=F4*B22*(10^-18)*B14"
$ws.Range("K5").Formula = "=F5*B22*(10^-18)*B14"
$ws.Range("K6").Formula = "=F6*B22*(10^-18)*B14"
$ws.Range("K7").Formula = "=F7*B22*(10^-18)*B14"

# K17 is the trailing header-row cell above column K, left blank but styled
# like its neighbours (J17, I17, ...). Copy the format from J17.
$ws.Range("J17").Copy()
$ws.Range("K17").PasteSpecial(-4122)

# Column K width, matching the new column's authored width.
$ws.Columns.Item(11).ColumnWidth = 15.02

# Update the remembered selection to match the authored state.
$ws.Range("L9").Select()
